# [Typed SDK] Fix problems found during testing Explorer example
#
# 1. Rename "Sheet1" -> "List of classes"
# 2. Fix a typo in the Events column of the "DataBrowser" row (G12):
#    "Many manye events" -> "Many many events"
# 3. Make the "List of classes" sheet the active tab/sheet, with cell
#    G13 selected (instead of "Intro" being active with A63 selected).

$wb = $excel.ActiveWorkbook

$wsClasses = $wb.Worksheets.Item("Sheet1")
$wsClasses.Name = "List of classes"

$wsClasses.Range("G12").Value = "Many many events"

$wsClasses.Activate()
$wsClasses.Range("G13").Select()
